# Adds newly-uploaded repo rows to the "software_tools", "tutorials" and
# "manuscripts" sheets, then leaves the selection/active-sheet state the
# way the author left it when they saved (software_tools active/selected,
# with the cursor parked at I25; manuscripts cursor at C2; tutorials cursor
# at D20; analysis cursor at A35, no longer the active tab).

$wb = $excel.ActiveWorkbook

$tools       = $wb.Worksheets.Item("software_tools")
$tutorials   = $wb.Worksheets.Item("tutorials")
$manuscripts = $wb.Worksheets.Item("manuscripts")
$analysis    = $wb.Worksheets.Item("analysis")

# ---------------------------------------------------------------------
# software_tools (sheet1) - rows 11-17
# ---------------------------------------------------------------------
$tools.Cells.Item(11, 1).Value = "hsvfinder"
$tools.Cells.Item(11, 2).Value = "Magic Wand implementation of hsv region finding in images"
$tools.Cells.Item(11, 3).Value = "https://github.com/TeamMacLean/hsvfinder"
$tools.Cells.Item(11, 4).Value = "danmaclean"

$tools.Cells.Item(12, 1).Value = "besthr"
$tools.Cells.Item(12, 2).Value = "bootstrap estimation of hr scores"
$tools.Cells.Item(12, 3).Value = "https://github.com/TeamMacLean/besthr"
$tools.Cells.Item(12, 4).Value = "danmaclean"

$tools.Cells.Item(13, 1).Value = "fiat"
$tools.Cells.Item(13, 2).Value = "general fluorescence image analysis"
$tools.Cells.Item(13, 3).Value = "https://github.com/TeamMacLean/fiat"
$tools.Cells.Item(13, 4).Value = "danmaclean"

$tools.Cells.Item(14, 1).Value = "glowseeds"
$tools.Cells.Item(14, 2).Value = "count fluorescent seeds in bf and fluo images"
$tools.Cells.Item(14, 3).Value = "https://github.com/TeamMacLean/glowseeds"
$tools.Cells.Item(14, 4).Value = "danmaclean"

# row 15 was typed link-first, then description, then name
$tools.Cells.Item(15, 3).Value = "https://github.com/TeamMacLean/stomatameasurer"
$tools.Cells.Item(15, 2).Value = "count stomata in FLEX images"
$tools.Cells.Item(15, 1).Value = "stomatameasurer"
$tools.Cells.Item(15, 4).Value = "danmaclean"

$tools.Cells.Item(16, 1).Value = "redpatchcam"
$tools.Cells.Item(16, 2).Value = "GUI for redpatch on Raspberry Pi"
$tools.Cells.Item(16, 3).Value = "https://github.com/TeamMacLean/redpatchcam"
$tools.Cells.Item(16, 4).Value = "danmaclean"

$tools.Cells.Item(17, 1).Value = "fluorseg"
$tools.Cells.Item(17, 2).Value = "Segment fluorescent images"
$tools.Cells.Item(17, 3).Value = "https://github.com/TeamMacLean/fluorseg"
$tools.Cells.Item(17, 4).Value = "danmaclean"

# ---------------------------------------------------------------------
# tutorials (sheet4) - row 7 ("basic_alignment")
# ---------------------------------------------------------------------
$tutorials.Cells.Item(7, 1).Value = "basic_alignment"
$tutorials.Cells.Item(7, 2).Value = "Tutorial on basic alignment on HPC"
$tutorials.Cells.Item(7, 3).Value = "https://github.com/TeamMacLean/basic_alignment"
$tutorials.Cells.Item(7, 4).Value = "danmaclean"

# ---------------------------------------------------------------------
# software_tools (sheet1) - rows 18-23
# ---------------------------------------------------------------------
$tools.Cells.Item(18, 1).Value = "phobiius_galaxy"
$tools.Cells.Item(18, 2).Value = "Wrapper to run phobius in Galaxy"
$tools.Cells.Item(18, 3).Value = "https://github.com/TeamMacLean/phobius_wrapper_tool"
$tools.Cells.Item(18, 4).Value = "tsl-ramkrishna"

$tools.Cells.Item(19, 1).Value = "signalp_galaxy_wrapper"
$tools.Cells.Item(19, 2).Value = "Wrapper to run SignalP in galaxy"
$tools.Cells.Item(19, 3).Value = "https://github.com/TeamMacLean/signalp_galaxy_wrappers"
$tools.Cells.Item(19, 4).Value = "tsl-ramkrishna"

$tools.Cells.Item(20, 1).Value = "kamoun_image_analysis_tools"
$tools.Cells.Item(20, 2).Value = "Some ia tools for kamoun group"
$tools.Cells.Item(20, 3).Value = "https://github.com/TeamMacLean/kamoun_image_analysis_tools"
$tools.Cells.Item(20, 4).Value = "danmaclean"

$tools.Cells.Item(21, 1).Value = "atacr"
$tools.Cells.Item(21, 2).Value = "Analysin cap seq count data"
$tools.Cells.Item(21, 3).Value = "https://github.com/TeamMacLean/atacr"
$tools.Cells.Item(21, 4).Value = "danmaclean"

$tools.Cells.Item(22, 1).Value = "peak_Caller"
$tools.Cells.Item(22, 2).Value = "finding peaks on waves from 384 well plates"
$tools.Cells.Item(22, 3).Value = "https://github.com/TeamMacLean/peak_caller"
$tools.Cells.Item(22, 4).Value = "danmaclean"

$tools.Cells.Item(23, 1).Value = "geefu.io"
$tools.Cells.Item(23, 2).Value = "JS version of genome browser"
$tools.Cells.Item(23, 3).Value = "https://github.com/TeamMacLean/geefu.io"
$tools.Cells.Item(23, 4).Value = "martinpage"

# ---------------------------------------------------------------------
# tutorials (sheet4) - row 8 ("ss_2017")
# ---------------------------------------------------------------------
$tutorials.Cells.Item(8, 1).Value = "ss_2017"
$tutorials.Cells.Item(8, 2).Value = "Summer School 2017 Handbook"
$tutorials.Cells.Item(8, 3).Value = "https://github.com/TeamMacLean/summer_school_handbook"
$tutorials.Cells.Item(8, 4).Value = "danmaclean"

# ---------------------------------------------------------------------
# manuscripts (sheet3) - row 2
# ---------------------------------------------------------------------
$manuscripts.Cells.Item(2, 1).Value = "pilar_corredor_lab_book"
$manuscripts.Cells.Item(2, 2).Value = "https://github.com/TeamMacLean/Lab_book_TSL"

# ---------------------------------------------------------------------
# software_tools (sheet1) - rows 24-28
# ---------------------------------------------------------------------
$tools.Cells.Item(24, 1).Value = "pdist"
$tools.Cells.Item(24, 2).Value = "ruby distance methods"
$tools.Cells.Item(24, 3).Value = "https://github.com/TeamMacLean/pdist"
$tools.Cells.Item(24, 4).Value = "edwardchalstrey"

$tools.Cells.Item(25, 1).Value = "candisnp"
$tools.Cells.Item(25, 2).Value = "Find candidate SNPs webapp"
$tools.Cells.Item(25, 3).Value = "https://github.com/TeamMacLean/candisnp"
$tools.Cells.Item(25, 4).Value = "danmaclean"

$tools.Cells.Item(26, 1).Value = "bioruby-samtools"
$tools.Cells.Item(26, 2).Value = "samtools in ruby"
$tools.Cells.Item(26, 3).Value = "https://github.com/TeamMacLean/bioruby-samtools"
$tools.Cells.Item(26, 4).Value = "danmaclean"

$tools.Cells.Item(27, 1).Value = "bioruby-svgenes"
$tools.Cells.Item(27, 2).Value = "render genes in ruby"
$tools.Cells.Item(27, 3).Value = "https://github.com/TeamMacLean/bioruby-svgenes"
$tools.Cells.Item(27, 4).Value = "danmaclean"

$tools.Cells.Item(28, 1).Value = "geefu"
$tools.Cells.Item(28, 2).Value = "genome browser"
$tools.Cells.Item(28, 3).Value = "https://github.com/TeamMacLean/gee_fu"
$tools.Cells.Item(28, 4).Value = "danmaclean"

# ---------------------------------------------------------------------
# Leave the selections/active sheet matching the saved file: the
# manuscripts and tutorials sheets keep the cursor where the author left
# it, the previously-active "analysis" sheet is parked elsewhere and is
# no longer the active tab, and "software_tools" ends up active/selected.
# ---------------------------------------------------------------------
$manuscripts.Range("C2").Select() | Out-Null
$tutorials.Range("D20").Select() | Out-Null
$analysis.Range("A35").Select() | Out-Null

$tools.Activate()
$tools.Range("I25").Select() | Out-Null
